$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 44, shifting existing rows 44-101 down to 45-102.
$ws.Rows.Item(44).Insert()

# Populate the newly inserted row 44 with the new weekly record.
$ws.Cells.Item(44, 1).Value = 7
$ws.Cells.Item(44, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(44, 3).Value = "Ñuble"
$ws.Cells.Item(44, 4).Value = 45125
$ws.Cells.Item(44, 5).Value = 16
$ws.Cells.Item(44, 6).Value = 100112044
$ws.Cells.Item(44, 7).Value = "Perejil"
$ws.Cells.Item(44, 8).Value = "Sin especificar"
$ws.Cells.Item(44, 9).Value = "Primera"
$ws.Cells.Item(44, 10).Value = 100
$ws.Cells.Item(44, 11).Value = 1500
$ws.Cells.Item(44, 12).Value = 1500
$ws.Cells.Item(44, 13).Value = 1500
$ws.Cells.Item(44, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(44, 15).Value = "Región de Ñuble"
$ws.Cells.Item(44, 16).Value = 1500
$ws.Cells.Item(44, 17).Value = 1
$ws.Cells.Item(44, 18).Value = "Hortaliza"
